# "xlswrite - filtros con lista"
#
# Adds a second worksheet ("Ficha con filtros de lista") next to the
# existing "Ficha con filtros por condición" sheet. The new sheet lists
# the 12 months with their month-number, and has an AutoFilter on column
# A restricted to a fixed list of values (4, 8, 11 -> Abril, Agosto,
# Noviembre), hiding the rows that don't match.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Keep a plain (non workbook-qualified) reference to sheet1's existing
# filtered range before we touch anything, so we can recreate its hidden
# _FilterDatabase name further down in the order the target file uses.
$sheet1Ref = "'" + $ws1.Name + "'!" + $ws1.Range("A1:D10").Address()

# Drop any pre-existing _FilterDatabase defined names - they'll be
# recreated below (new sheet first, then the original sheet) so the
# <definedNames> order matches the target workbook.
for ($i = $wb.Names.Count; $i -ge 1; $i--) {
    $nm = $wb.Names.Item($i)
    if ($nm.Name.EndsWith("_FilterDatabase")) {
        $nm.Delete()
    }
}

# New sheet, placed right after the first one.
$ws2 = $wb.Worksheets.Add([Type]::Missing, $ws1)
$ws2.Name = "Ficha con filtros de lista"

# Header row.
$ws2.Range("A1").Value = "Nro. Mes"
$ws2.Range("B1").Value = "Mes"
$ws2.Range("A1:B1").Font.Bold = $true

# 12 data rows: month number + month name.
$meses = @("Enero","Febrero","Marzo","Abril","Mayo","Junio","Julio","Agosto","Septiembre","Octubre","Noviembre","Diciembre")
for ($i = 0; $i -lt 12; $i++) {
    $row = $i + 2
    $ws2.Cells.Item($row, 1).Value = $i + 1
    $ws2.Cells.Item($row, 2).Value = $meses[$i]
}

# Column widths (approximate - narrow "Nro. Mes" / "Mes" columns).
$ws2.Columns.Item(1).ColumnWidth = 5.75
$ws2.Columns.Item(2).ColumnWidth = 9.75

# AutoFilter on column A (Nro. Mes) restricted to a fixed value list:
# 4 (Abril), 8 (Agosto), 11 (Noviembre). xlFilterValues = 7. This hides
# every other row automatically, same as Excel's "filtro de lista".
$ws2.Range("A1:B13").AutoFilter(1, @("4","8","11"), 7) | Out-Null

# Recreate the hidden _xlnm._FilterDatabase names - new sheet's first,
# then the original sheet's - matching the target file's order.
$sheet2Ref = "'" + $ws2.Name + "'!" + $ws2.Range("A1:B13").Address()
$ws2.Names.Add("_xlnm._FilterDatabase", "=" + $sheet2Ref, $false) | Out-Null
$ws1.Names.Add("_xlnm._FilterDatabase", "=" + $sheet1Ref, $false) | Out-Null
$wb.Names.Item(1).Visible = $false
$wb.Names.Item(2).Visible = $false

# Leave the original sheet as the active/selected one, like before the edit.
$ws1.Activate()
$ws1.Select()
